# ContosoLearn Market Research — translate Indonesian bullet text to English
# and merge the multi-run paragraphs into single runs, matching the target
# OOXML. The WoodgroveLMS bullet additionally gets a grammar-check style
# "a best" flag, expressed with <w:proofErr w:type="gramStart"/>/"gramEnd"
# markers splitting that bullet into three runs.

$d = $word.ActiveDocument

function Set-ParagraphText($paraIndex, $newText) {
    $p = $d.Paragraphs($paraIndex)
    $r = $p.Range
    # Exclude the trailing paragraph mark from the range we overwrite.
    $r.End = $r.End - 1
    $r.Text = $newText
}

# 1. Title
Set-ParagraphText 1 "ContosoLearn Market Research"

# 2. AdatumLearn bullet (merges 3 runs into 1)
Set-ParagraphText 2 "AdatumLearn: AdatumLearn is a top AI-powered learning platform that uses artificial intelligence to enrich eLearning with features that automate a variety of tasks. It is known for its content authoring capabilities and adaptive learning technology."

# 3. AdventureLearn bullet
Set-ParagraphText 3 "AdventureLearn: AdventureLearn is another AI-powered learning platform that offers personalized learning experiences and data-driven recommendations."

# 4. AlpineTraining bullet
Set-ParagraphText 4 "AlpineTraining: AlpineTraining is a mobile-first learning platform that focuses on microlearning."

# 5. Bellows OnDemand bullet
Set-ParagraphText 5 "Bellows OnDemand: Bellows OnDemand is a comprehensive learning solution that offers content creation and social collaboration."

# 6. FabrikamLearning bullet
Set-ParagraphText 6 "FabrikamLearning: FabrikamLearning provides a suite of learning platforms that cater to different learning needs."

# 7. FirstUp Cards bullet
Set-ParagraphText 7 "FirstUp Cards: FirstUp Cards is a mobile learning app that is ideal for training on safety procedures, compliance, new product knowledge or any other type of training scenario."

# 8. Munson'sLearn bullet
Set-ParagraphText 8 "Munson'sLearn: Munson'sLearn is designed to enable businesses to train their employees, partners, and customers."

# 9. LibertyLearn bullet
Set-ParagraphText 9 "LibertyLearn: LibertyLearn is a fast LMS for your mission-critical project."

# 10. WoodgroveLMS bullet — needs to end up as three runs split around "a best",
# separated by <w:proofErr w:type="gramStart"/> ... <w:proofErr w:type="gramEnd"/>,
# mimicking Word's grammar-checker markup. Plain Range.Text cannot create those
# proofErr siblings, so build the paragraph from raw WordOpenXML via InsertXML,
# while preserving the original pPr (ListParagraph / numPr) of the bullet.
$p10 = $d.Paragraphs(10)
$r10 = $p10.Range
$r10.End = $r10.End - 1
$r10.Text = ""
$r10b = $d.Paragraphs(10).Range
$r10b.End = $r10b.End - 1
$xml10 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + `
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
  '<pkg:xmlData>' + `
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
  '<w:body>' + `
  '<w:p>' + `
  '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="14"/></w:numPr></w:pPr>' + `
  '<w:r><w:t xml:space="preserve">WoodgroveLMS: WoodgroveLMS is a functional and attractive learning management system built to provide </w:t></w:r>' + `
  '<w:proofErr w:type="gramStart"/>' + `
  '<w:r><w:t>a best</w:t></w:r>' + `
  '<w:proofErr w:type="gramEnd"/>' + `
  '<w:r><w:t>-in-class training experience.</w:t></w:r>' + `
  '</w:p>' + `
  '</w:body></w:document>' + `
  '</pkg:xmlData></pkg:part></pkg:package>'
[void]$r10b.InsertXML($xml10)

# 11. NorthwindWorlds bullet
Set-ParagraphText 11 "NorthwindWorlds: NorthwindWorlds is a powerful, easy-to-use, and reliable training solution for individuals and enterprises."

# 12. ProsewareLearn bullet
Set-ParagraphText 12 "ProsewareLearn: ProsewareLearn is an online education company that offers a variety of video training courses for software developers, IT administrators, and creative professionals through its website."

# 13. RelecloudLearn bullet
Set-ParagraphText 13 "RelecloudLearn: RelecloudLearn is an American online learning platform that offers massive open online courses (MOOC), specializations, and degrees in a variety of subjects."

# 14. TreyAcademy bullet
Set-ParagraphText 14 "TreyAcademy: TreyAcademy is an online learning platform aimed at professional adults and students, developed in May 2010."

# 15. Closing paragraph (merges 4 runs into 1)
Set-ParagraphText 15 "These platforms have a significant market presence and are widely recognized for their AI-powered features, such as personalized learning experiences, data-driven recommendations, and automation of tasks. They are transforming the eLearning landscape by leveraging AI to deliver more engaging, rewarding, and personalized learning experiences. "
